$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.935093
$ws.Range("H2").Value = 98.80527900000001
$ws.Range("I2").Value = 0.247867546118624
$ws.Range("J2").Value = 0.247867546118624
$ws.Range("M2").Value = 6.970140000000001
$ws.Range("N2").Value = 20.91042
$ws.Range("O2").Value = 0.06638545261649673
$ws.Range("P2").Value = 0.06638545261649674
$ws.Range("Q2").Value = 229.56220912302
$ws.Range("R2").Value = 2066.05988210718
$ws.Range("S2").Value = 0.01645479923802523
$ws.Range("T2").Value = 0.01645479923802523
$ws.Range("G3").Value = 32.935093
$ws.Range("H3").Value = 98.80527900000001
$ws.Range("I3").Value = 0.247867546118624
$ws.Range("J3").Value = 0.247867546118624
$ws.Range("O3").Value = 0.03236297878883257
$ws.Range("P3").Value = 0.03236297878883258
$ws.Range("Q3").Value = 111.9118212160787
$ws.Range("R3").Value = 1007.206390944708
$ws.Range("S3").Value = 0.008021732137477007
$ws.Range("T3").Value = 0.008021732137477009
$ws.Range("G4").Value = 32.935093
$ws.Range("H4").Value = 98.80527900000001
$ws.Range("I4").Value = 0.247867546118624
$ws.Range("J4").Value = 0.247867546118624
$ws.Range("M4").Value = 40.232648
$ws.Range("N4").Value = 120.697944
$ws.Range("O4").Value = 0.3831863560043545
$ws.Range("P4").Value = 0.3831863560043545
$ws.Range("Q4").Value = 1325.066003516264
$ws.Range("R4").Value = 11925.59403164638
$ws.Range("S4").Value = 0.09497946176893682
$ws.Range("T4").Value = 0.09497946176893682
$ws.Range("G5").Value = 32.935093
$ws.Range("H5").Value = 98.80527900000001
$ws.Range("I5").Value = 0.247867546118624
$ws.Range("J5").Value = 0.247867546118624
$ws.Range("M5").Value = 0.6731889999999999
$ws.Range("N5").Value = 2.019567
$ws.Range("O5").Value = 0.006411629674790867
$ws.Range("P5").Value = 0.006411629674790868
$ws.Range("Q5").Value = 22.171542321577
$ws.Range("R5").Value = 199.543880894193
$ws.Range("S5").Value = 0.001589234914111763
$ws.Range("T5").Value = 0.001589234914111764
$ws.Range("G6").Value = 32.935093
$ws.Range("H6").Value = 98.80527900000001
$ws.Range("I6").Value = 0.247867546118624
$ws.Range("J6").Value = 0.247867546118624
$ws.Range("M6").Value = 53.72106333333333
$ws.Range("N6").Value = 161.16319
$ws.Range("O6").Value = 0.5116535829155252
$ws.Range("P6").Value = 0.5116535829155253
$ws.Range("Q6").Value = 1769.308216942223
$ws.Range("R6").Value = 15923.77395248001
$ws.Range("S6").Value = 0.1268223180600732
$ws.Range("T6").Value = 0.1268223180600732
$ws.Range("I7").Value = 0.4962147730988433
$ws.Range("J7").Value = 0.4962147730988432
$ws.Range("M7").Value = 6.970140000000001
$ws.Range("N7").Value = 20.91042
$ws.Range("O7").Value = 0.06638545261649673
$ws.Range("P7").Value = 0.06638545261649674
$ws.Range("Q7").Value = 459.56867405922
$ws.Range("R7").Value = 4136.11806653298
$ws.Range("S7").Value = 0.03294144230715893
$ws.Range("T7").Value = 0.03294144230715894
$ws.Range("I8").Value = 0.4962147730988433
$ws.Range("J8").Value = 0.4962147730988432
$ws.Range("O8").Value = 0.03236297878883257
$ws.Range("P8").Value = 0.03236297878883258
$ws.Range("S8").Value = 0.01605898817650323
$ws.Range("T8").Value = 0.01605898817650323
$ws.Range("I9").Value = 0.4962147730988433
$ws.Range("J9").Value = 0.4962147730988432
$ws.Range("M9").Value = 40.232648
$ws.Range("N9").Value = 120.697944
$ws.Range("O9").Value = 0.3831863560043545
$ws.Range("P9").Value = 0.3831863560043545
$ws.Range("Q9").Value = 2652.696315318104
$ws.Range("R9").Value = 23874.26683786293
$ws.Range("S9").Value = 0.1901427306992733
$ws.Range("T9").Value = 0.1901427306992733
$ws.Range("I10").Value = 0.4962147730988433
$ws.Range("J10").Value = 0.4962147730988432
$ws.Range("M10").Value = 0.6731889999999999
$ws.Range("N10").Value = 2.019567
$ws.Range("O10").Value = 0.006411629674790867
$ws.Range("P10").Value = 0.006411629674790868
$ws.Range("Q10").Value = 44.38599169044699
$ws.Range("R10").Value = 399.4739252140229
$ws.Range("S10").Value = 0.00318154536427016
$ws.Range("T10").Value = 0.00318154536427016
$ws.Range("I11").Value = 0.4962147730988433
$ws.Range("J11").Value = 0.4962147730988432
$ws.Range("M11").Value = 53.72106333333333
$ws.Range("N11").Value = 161.16319
$ws.Range("O11").Value = 0.5116535829155252
$ws.Range("P11").Value = 0.5116535829155253
$ws.Range("Q11").Value = 3542.040453298122
$ws.Range("R11").Value = 31878.3640796831
$ws.Range("S11").Value = 0.2538900665516375
$ws.Range("T11").Value = 0.2538900665516375
$ws.Range("G12").Value = 13.46314666666667
$ws.Range("H12").Value = 40.38944
$ws.Range("I12").Value = 0.1013228390550407
$ws.Range("J12").Value = 0.1013228390550407
$ws.Range("M12").Value = 6.970140000000001
$ws.Range("N12").Value = 20.91042
$ws.Range("O12").Value = 0.06638545261649673
$ws.Range("P12").Value = 0.06638545261649674
$ws.Range("Q12").Value = 93.84001710720001
$ws.Range("R12").Value = 844.5601539648001
$ws.Range("S12").Value = 0.006726362531057331
$ws.Range("T12").Value = 0.006726362531057332
$ws.Range("G13").Value = 13.46314666666667
$ws.Range("H13").Value = 40.38944
$ws.Range("I13").Value = 0.1013228390550407
$ws.Range("J13").Value = 0.1013228390550407
$ws.Range("O13").Value = 0.03236297878883257
$ws.Range("P13").Value = 0.03236297878883258
$ws.Range("Q13").Value = 45.74710819143111
$ws.Range("R13").Value = 411.72397372288
$ws.Range("S13").Value = 0.00327910889116258
$ws.Range("T13").Value = 0.00327910889116258
$ws.Range("G14").Value = 13.46314666666667
$ws.Range("H14").Value = 40.38944
$ws.Range("I14").Value = 0.1013228390550407
$ws.Range("J14").Value = 0.1013228390550407
$ws.Range("M14").Value = 40.232648
$ws.Range("N14").Value = 120.697944
$ws.Range("O14").Value = 0.3831863560043545
$ws.Range("P14").Value = 0.3831863560043545
$ws.Range("Q14").Value = 541.6580408123733
$ws.Range("R14").Value = 4874.92236731136
$ws.Range("S14").Value = 0.03882552947751675
$ws.Range("T14").Value = 0.03882552947751675
$ws.Range("G15").Value = 13.46314666666667
$ws.Range("H15").Value = 40.38944
$ws.Range("I15").Value = 0.1013228390550407
$ws.Range("J15").Value = 0.1013228390550407
$ws.Range("M15").Value = 0.6731889999999999
$ws.Range("N15").Value = 2.019567
$ws.Range("O15").Value = 0.006411629674790867
$ws.Range("P15").Value = 0.006411629674790868
$ws.Range("Q15").Value = 9.063242241386666
$ws.Range("R15").Value = 81.56918017248
$ws.Range("S15").Value = 0.0006496445216193582
$ws.Range("T15").Value = 0.0006496445216193583
$ws.Range("G16").Value = 13.46314666666667
$ws.Range("H16").Value = 40.38944
$ws.Range("I16").Value = 0.1013228390550407
$ws.Range("J16").Value = 0.1013228390550407
$ws.Range("M16").Value = 53.72106333333333
$ws.Range("N16").Value = 161.16319
$ws.Range("O16").Value = 0.5116535829155252
$ws.Range("P16").Value = 0.5116535829155253
$ws.Range("Q16").Value = 723.2545547459555
$ws.Range("R16").Value = 6509.2909927136
$ws.Range("S16").Value = 0.0518421936336847
$ws.Range("T16").Value = 0.05184219363368471
$ws.Range("G17").Value = 6.098311666666667
$ws.Range("H17").Value = 18.294935
$ws.Range("I17").Value = 0.04589552998326869
$ws.Range("J17").Value = 0.04589552998326869
$ws.Range("M17").Value = 6.970140000000001
$ws.Range("N17").Value = 20.91042
$ws.Range("O17").Value = 0.06638545261649673
$ws.Range("P17").Value = 0.06638545261649674
$ws.Range("Q17").Value = 42.50608608030001
$ws.Range("R17").Value = 382.5547747227001
$ws.Range("S17").Value = 0.003046795531013289
$ws.Range("T17").Value = 0.003046795531013289
$ws.Range("G18").Value = 6.098311666666667
$ws.Range("H18").Value = 18.294935
$ws.Range("I18").Value = 0.04589552998326869
$ws.Range("J18").Value = 0.04589552998326869
$ws.Range("O18").Value = 0.03236297878883257
$ws.Range("P18").Value = 0.03236297878883258
$ws.Range("Q18").Value = 20.72176219329111
$ws.Range("R18").Value = 186.49585973962
$ws.Range("S18").Value = 0.001485316063350754
$ws.Range("T18").Value = 0.001485316063350754
$ws.Range("G19").Value = 6.098311666666667
$ws.Range("H19").Value = 18.294935
$ws.Range("I19").Value = 0.04589552998326869
$ws.Range("J19").Value = 0.04589552998326869
$ws.Range("M19").Value = 40.232648
$ws.Range("N19").Value = 120.697944
$ws.Range("O19").Value = 0.3831863560043545
$ws.Range("P19").Value = 0.3831863560043545
$ws.Range("Q19").Value = 245.3512266792933
$ws.Range("R19").Value = 2208.16104011364
$ws.Range("S19").Value = 0.01758654089117732
$ws.Range("T19").Value = 0.01758654089117732
$ws.Range("G20").Value = 6.098311666666667
$ws.Range("H20").Value = 18.294935
$ws.Range("I20").Value = 0.04589552998326869
$ws.Range("J20").Value = 0.04589552998326869
$ws.Range("M20").Value = 0.6731889999999999
$ws.Range("N20").Value = 2.019567
$ws.Range("O20").Value = 0.006411629674790867
$ws.Range("P20").Value = 0.006411629674790868
$ws.Range("Q20").Value = 4.105316332571666
$ws.Range("R20").Value = 36.947846993145
$ws.Range("S20").Value = 0.0002942651419809795
$ws.Range("T20").Value = 0.0002942651419809796
$ws.Range("G21").Value = 6.098311666666667
$ws.Range("H21").Value = 18.294935
$ws.Range("I21").Value = 0.04589552998326869
$ws.Range("J21").Value = 0.04589552998326869
$ws.Range("M21").Value = 53.72106333333333
$ws.Range("N21").Value = 161.16319
$ws.Range("O21").Value = 0.5116535829155252
$ws.Range("P21").Value = 0.5116535829155253
$ws.Range("Q21").Value = 327.6077872714055
$ws.Range("R21").Value = 2948.47008544265
$ws.Range("S21").Value = 0.02348261235574634
$ws.Range("T21").Value = 0.02348261235574635
$ws.Range("G22").Value = 14.44328633333333
$ws.Range("H22").Value = 43.329859
$ws.Range("I22").Value = 0.1086993117442235
$ws.Range("J22").Value = 0.1086993117442234
$ws.Range("M22").Value = 6.970140000000001
$ws.Range("N22").Value = 20.91042
$ws.Range("O22").Value = 0.06638545261649673
$ws.Range("P22").Value = 0.06638545261649674
$ws.Range("Q22").Value = 100.67172780342
$ws.Range("R22").Value = 906.0455502307801
$ws.Range("S22").Value = 0.007216053009241952
$ws.Range("T22").Value = 0.007216053009241953
$ws.Range("G23").Value = 14.44328633333333
$ws.Range("H23").Value = 43.329859
$ws.Range("I23").Value = 0.1086993117442235
$ws.Range("J23").Value = 0.1086993117442234
$ws.Range("O23").Value = 0.03236297878883257
$ws.Range("P23").Value = 0.03236297878883258
$ws.Range("Q23").Value = 49.07757442520755
$ws.Range("R23").Value = 441.698169826868
$ws.Range("S23").Value = 0.003517833520339003
$ws.Range("T23").Value = 0.003517833520339003
$ws.Range("G24").Value = 14.44328633333333
$ws.Range("H24").Value = 43.329859
$ws.Range("I24").Value = 0.1086993117442235
$ws.Range("J24").Value = 0.1086993117442234
$ws.Range("M24").Value = 40.232648
$ws.Range("N24").Value = 120.697944
$ws.Range("O24").Value = 0.3831863560043545
$ws.Range("P24").Value = 0.3831863560043545
$ws.Range("Q24").Value = 581.0916550122106
$ws.Range("R24").Value = 5229.824895109896
$ws.Range("S24").Value = 0.04165209316745033
$ws.Range("T24").Value = 0.04165209316745032
$ws.Range("G25").Value = 14.44328633333333
$ws.Range("H25").Value = 43.329859
$ws.Range("I25").Value = 0.1086993117442235
$ws.Range("J25").Value = 0.1086993117442234
$ws.Range("M25").Value = 0.6731889999999999
$ws.Range("N25").Value = 2.019567
$ws.Range("O25").Value = 0.006411629674790867
$ws.Range("P25").Value = 0.006411629674790868
$ws.Range("Q25").Value = 9.723061483450332
$ws.Range("R25").Value = 87.50755335105299
$ws.Range("S25").Value = 0.0006969397328086065
$ws.Range("T25").Value = 0.0006969397328086065
$ws.Range("G26").Value = 14.44328633333333
$ws.Range("H26").Value = 43.329859
$ws.Range("I26").Value = 0.1086993117442235
$ws.Range("J26").Value = 0.1086993117442234
$ws.Range("M26").Value = 53.72106333333333
$ws.Range("N26").Value = 161.16319
$ws.Range("O26").Value = 0.5116535829155252
$ws.Range("P26").Value = 0.5116535829155253
$ws.Range("Q26").Value = 775.9086998544676
$ws.Range("R26").Value = 6983.178298690209
$ws.Range("S26").Value = 0.05561639231438355
$ws.Range("T26").Value = 0.05561639231438356
